$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text updates (Volume number, week-of dates) ----
$ws.Range("A8").Value = "Volume 32   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/15/2025  Through  12/21/2025"

# ---- Weekly crime-stat table (rows 15-28, 31) ----

# -- Plain numeric value updates (style unchanged) --
$ws.Range("N15").Value = -56.756756756756
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -30
$ws.Range("I16").Value = 112
$ws.Range("J16").Value = 139
$ws.Range("K16").Value = -19.424460431654
$ws.Range("L16").Value = -39.459459459459
$ws.Range("M16").Value = -40.425531914893
$ws.Range("N16").Value = -88.477366255144
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -63.636363636363
$ws.Range("I17").Value = 177
$ws.Range("J17").Value = 204
$ws.Range("K17").Value = -13.235294117647
$ws.Range("L17").Value = -19.17808219178
$ws.Range("M17").Value = 20.408163265306
$ws.Range("N17").Value = -68.892794376098
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -35.294117647058
$ws.Range("I18").Value = 164
$ws.Range("J18").Value = 172
$ws.Range("K18").Value = -4.651162790697
$ws.Range("L18").Value = -39.483394833948
$ws.Range("M18").Value = -23.00469483568
$ws.Range("N18").Value = -82.173913043478
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -6.666666666666
$ws.Range("F19").Value = 64
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = 12.280701754386
$ws.Range("I19").Value = 803
$ws.Range("J19").Value = 745
$ws.Range("K19").Value = 7.785234899328
$ws.Range("L19").Value = -17.045454545454
$ws.Range("M19").Value = 7.066666666666
$ws.Range("N19").Value = -51.568154402895
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 36
$ws.Range("K20").Value = -20
$ws.Range("L20").Value = -20
$ws.Range("M20").Value = -23.404255319148
$ws.Range("N20").Value = -93.296089385474
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 92
$ws.Range("G21").Value = 109
$ws.Range("H21").Value = -15.596330275229
$ws.Range("I21").Value = 1309
$ws.Range("J21").Value = 1325
$ws.Range("K21").Value = -1.207547169811
$ws.Range("L21").Value = -23.31575864089
$ws.Range("M21").Value = -3.961848862802
$ws.Range("N21").Value = -72.208067940552
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 18
$ws.Range("K22").Value = 38.461538461538
$ws.Range("L22").Value = 50
$ws.Range("M22").Value = -18.181818181818
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -75
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = -54.545454545454
$ws.Range("I23").Value = 122
$ws.Range("J23").Value = 95
$ws.Range("K23").Value = 28.421052631578
$ws.Range("L23").Value = 0.826446280991
$ws.Range("M23").Value = 11.926605504587
$ws.Range("C24").Value = 41
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = 41.379310344827
$ws.Range("F24").Value = 143
$ws.Range("G24").Value = 117
$ws.Range("H24").Value = 22.222222222222
$ws.Range("I24").Value = 1574
$ws.Range("J24").Value = 1533
$ws.Range("K24").Value = 2.674494455316
$ws.Range("L24").Value = 6.784260515603
$ws.Range("M24").Value = -7.737397420867
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = 28.571428571428
$ws.Range("F25").Value = 68
$ws.Range("G25").Value = 52
$ws.Range("H25").Value = 30.76923076923
$ws.Range("I25").Value = 797
$ws.Range("J25").Value = 873
$ws.Range("K25").Value = -8.705612829324
$ws.Range("L25").Value = -3.510895883777
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -12.5
$ws.Range("F26").Value = 37
$ws.Range("G26").Value = 36
$ws.Range("H26").Value = 2.777777777777
$ws.Range("I26").Value = 447
$ws.Range("J26").Value = 433
$ws.Range("K26").Value = 3.233256351039
$ws.Range("L26").Value = -3.663793103448
$ws.Range("M26").Value = -2.826086956521
$ws.Range("H27").Value = -100
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = -20
$ws.Range("I28").Value = 46
$ws.Range("J28").Value = 66
$ws.Range("K28").Value = -30.30303030303
$ws.Range("L28").Value = -25.806451612903

# -- Cells changing FROM text-placeholder style TO numeric style (need format copied from a numeric-style cell) --
$ws.Range("C16").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 2

# -- Cells changing FROM numeric style TO text-placeholder ("0" / "***.*") --
# Step 1: force text storage via Text number format, write placeholder string
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "0"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "***.*"

# Step 2: restore the General/style-13 formatting (reuses existing style instead of the temp Text style)
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("H31").PasteSpecial(-4122)

$excel.CutCopyMode = 0
